$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so numeric-looking strings
# (e.g. "95.875.01", "1.00", "0.0000194") are preserved exactly as authored,
# instead of being auto-converted into numbers by Excel.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = '95.875.01'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '3.309.87'
$ws.Range("E3").Value = '  -3.43%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '247.60'
$ws.Range("E5").Value = '  -3.85%  '
$ws.Range("D6").Value = '647.04'
$ws.Range("E6").Value = '  -1.74%  '
$ws.Range("D7").Value = '1.34'
$ws.Range("E7").Value = '  -10.62%  '
$ws.Range("D8").Value = '0.412'
$ws.Range("E8").Value = '  -4.53%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '0.967'
$ws.Range("E10").Value = '  -8.54%  '
$ws.Range("D11").Value = '3.306.62'
$ws.Range("E11").Value = '  -3.33%  '
$ws.Range("E12").Value = '  -4.63%  '
$ws.Range("D13").Value = '39.56'
$ws.Range("E13").Value = '  -5.97%  '
$ws.Range("D14").Value = '95.611.89'
$ws.Range("E14").Value = '  -2.16%  '
$ws.Range("D15").Value = '6.07'
$ws.Range("E15").Value = '  -6.04%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.0000248'
$ws.Range("E16").Value = '  -4.52%  '
$ws.Range("B17").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C17").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D17").Value = '3.924.97'
$ws.Range("E17").Value = '  -3.55%  '
$ws.Range("D18").Value = '8.47'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '3.316.91'
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").Value = '16.73'
$ws.Range("E20").Value = '  -5.00%  '
$ws.Range("D21").Value = '0.481'
$ws.Range("E21").Value = '  -6.76%  '
$ws.Range("D22").Value = '499.89'
$ws.Range("E22").Value = '  -2.05%  '
$ws.Range("D23").Value = '10.39'
$ws.Range("E23").Value = '  -5.94%  '
$ws.Range("D24").Value = '3.29'
$ws.Range("E24").Value = '  -5.24%  '
$ws.Range("D25").Value = '0.0000194'
$ws.Range("E25").Value = '  -6.33%  '
$ws.Range("D26").Value = '6.36'
$ws.Range("E26").Value = '  +4.49%  '
$ws.Range("D27").Value = '91.98'
$ws.Range("E27").Value = '  -4.68%  '
$ws.Range("D28").Value = '11.88'
$ws.Range("E28").Value = '  -7.36%  '
$ws.Range("D29").Value = '3.486.50'
$ws.Range("E29").Value = '  -3.40%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '0.141'
$ws.Range("E31").Value = '  -8.85%  '
$ws.Range("D32").Value = '10.73'
$ws.Range("E32").Value = '  -6.29%  '
$ws.Range("D33").Value = '0.183'
$ws.Range("E33").Value = '  -7.94%  '
$ws.Range("D34").Value = '2.44'
$ws.Range("E34").Value = '  +8.87%  '
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = '0.538'
$ws.Range("E36").Value = '  -6.45%  '
$ws.Range("D37").Value = '27.66'
$ws.Range("E37").Value = '  -7.31%  '
$ws.Range("D38").Value = '1.47'
$ws.Range("E38").Value = '  +3.92%  '
$ws.Range("D39").Value = '7.42'
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("E41").Value = '  -5.31%  '
$ws.Range("D42").Value = '500.45'
$ws.Range("E42").Value = '  -3.31%  '
$ws.Range("D43").Value = '24.28'
$ws.Range("E44").Value = '  -1.78%  '
$ws.Range("D45").Value = '0.817'
$ws.Range("E45").Value = '  -4.92%  '
$ws.Range("D46").Value = '0.0404'
$ws.Range("E46").Value = '  -4.14%  '
$ws.Range("D47").Value = '5.40'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").Value = '8.24'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("E49").Value = '  +0.68%  '
$ws.Range("D50").Value = '52.83'
$ws.Range("E50").Value = '  +3.23%  '
$ws.Range("D51").Value = '3.09'
$ws.Range("E51").Value = '  -7.22%  '

# Clean up: drop the temporary text-format override so the cells fall back
# to the workbook default style (no visible formatting change).
$priceVolumeRange.Style = "Normal"
